# The source workbook's "Item_2026" (column E) / "Monto_2026" label column (G) values
# were produced by an OCR pass that mangled many organism names (stray/garbled
# characters such as "TRil:lUTARIOS", "TECNOLOGiA", trailing boilerplate text
# copy-pasted from the PDF, etc.) and duplicated them as new shared strings
# instead of reusing the already-clean text that the sheet's E column carries
# for the same row. This script rewrites each mangled G-column label (and the
# one truncated E-column label) with the correct, clean text so that Excel's
# string de-duplication collapses them back onto a single shared string per
# organism, removing the stray OCR-garbled duplicates from the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E71 itself was truncated ("...TRÁNSITO Y") - restore the full organism name.
$ws.Range("E71").Value = 'AGENCIA NACIONAL DE TRÁNSITO Y SEGURIDAD VIAL'

# Column G (Monto_2026's paired label) fixes - one per affected row.
$ws.Range("G25").Value = 'JUSTICIA ELECTORAL'
$ws.Range("G30").Value = 'SINDICATURA GENERAL DE QUIEBRAS'
$ws.Range("G32").Value = 'DEFENSORÍA DEL PUEBLO'
$ws.Range("G35").Value = 'GOBIERNO DEPARTAMENTAL DE CONCEPCIÓN'
$ws.Range("G36").Value = 'GOBIERNO DEPARTAMENTAL DE SAN PEDRO'
$ws.Range("G37").Value = 'GOBIERNO DEPARTAMENTAL DE CORDILLERA'
$ws.Range("G43").Value = 'GOBIERNO DEPARTAMENTAL DE PARAGUARÍ'
$ws.Range("G44").Value = 'GOBIERNO DEPARTAMENTAL DE ALTO PARANÁ'
$ws.Range("G47").Value = 'GOBIERNO DEPARTAMENTAL DE AMAMBAY'
$ws.Range("G49").Value = 'GOBIERNO DEPARTAMENTAL DE PRESIDENTE HAYES'
$ws.Range("G52").Value = 'INSTITUTO NACIONAL DE TECNOLOGÍA, NORMALIZACIÓN Y METROLOGÍA'
$ws.Range("G55").Value = 'INSTITUTO PARAGUAYO DEL INDÍGENA'
$ws.Range("G66").Value = 'INSTITUTO PARAGUAYO DE TECNOLOGIA AGRARIA'
$ws.Range("G69").Value = 'SECRETARÍA DE DEFENSA DEL CONSUMIDOR Y EL USUARIO'
$ws.Range("G71").Value = 'AGENCIA NACIONAL DE TRÁNSITO Y SEGURIDAD VIAL'
$ws.Range("G73").Value = 'AGENCIA NACIONAL DE EVAL. Y ACRED. DE LA EDUCACIÓN SUPERIOR'
$ws.Range("G78").Value = 'DIRECCIÓN NACIONAL DE DEFENSA, SALUD Y BIENESTAR ANIMAL'
$ws.Range("G83").Value = 'DIRECCIÓN NACIONAL DE INGRESOS TRIBUTARIOS'
$ws.Range("G84").Value = 'INSTITUTO SUPERIOR NACIONAL DE MÚSICA'
$ws.Range("G90").Value = 'ADMINISTRACION NACIONAL DE ELECTRICIDAD'
$ws.Range("G100").Value = 'UNIVERSIDAD NACIONAL DEL ESTE'
$ws.Range("G105").Value = 'UNIVERSIDAD NACIONAL DE CAAGUAZÚ'
$ws.Range("G108").Value = 'UNIVERSIDAD NACIONAL DE MISIONES'
